$d = $word.ActiveDocument

# The API-reference table has two rows whose endpoint text ends in the same
# "?userId={id}&animeId={id}" query string:
#   ./danh-sach-yeu-thich/search/existsByAnime_IdAndUser_Id?userId={id}&animeId={id}
#   ./danh-sach-yeu-thich/search/countAllByAnime_IdAndUser_Id?userId={id}&animeId={id}
# Only the second ("countAllBy...") row is being updated, so the Find/Replace
# below uses the method name + query string together as the search target.
# That combined string is unique across the whole document, which guarantees
# the sibling "existsByAnime_IdAndUser_Id" row is left completely untouched.
$d.Content.Find.Execute(
    "countAllByAnime_IdAndUser_Id?userId={id}&animeId={id}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "countAllByUser_Id?userId={id}", 2) | Out-Null
